$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.174965858459473
$ws.Range("B1").Value = 2.352607488632202
$ws.Range("C1").Value = 3.408670663833618
$ws.Range("D1").Value = 1.679078459739685
$ws.Range("E1").Value = 1.213903665542603
